# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '43.121.27'
$ws.Range('E2').Value = '  -0.21%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '2.306.37'
$ws.Range('E3').Value = '  +0.04%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  +0.05%  '

# Row 5: BNB
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '300.80'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.42%  '

# Row 6: Solana
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '97.90'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.15%  '

# Row 7: XRP
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.518'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +2.60%  '

# Row 8: USDC
$ws.Range('E8').Value = '  +0.01%  '

# Row 9: Cardano
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.517'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.22%  '

# Row 10: Avalanche
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.85'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.18%  '

# Row 11: Dogecoin
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0792'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.25%  '

# Row 12: TRON
$ws.Range('E12').Value = '  +0.29%  '

# Row 13: Chainlink
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '17.96'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -4.78%  '

# Row 14: Polkadot
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.88'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.63%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range('D15').Value = '2.664.77'
$ws.Range('E15').Value = '  +0.00%  '

# Row 16: WrappedEther
$ws.Range('D16').Value = '2.297.01'
$ws.Range('E16').Value = '  +1.45%  '

# Row 17: Polygon
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.790'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.82%  '

# Row 18: WrappedBTC
$ws.Range('D18').Value = '43.023.91'
$ws.Range('E18').Value = '  -0.08%  '

# Row 19: InternetComputer(DFINITY)
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.14'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +3.47%  '

# Row 20: ShibaInu
$ws.Range('D20').Value = '0.0₃0911'
$ws.Range('E20').Value = '  +0.24%  '

# Row 21: Uniswap
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.13'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.31%  '

# Row 22: Litecoin
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.34'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.33%  '

# Row 23: BitcoinCash
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '238.10'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.53%  '

# Row 24: ImmutableX
$ws.Range('E24').Value = '  -1.67%  '

# Row 25: Dai
$ws.Range('E25').Value = '  -1.28%  '

# Row 26: LEO
$ws.Range('E26').Value = '  -0.25%  '

# Row 27: PancakeSwap
$ws.Range('E27').Value = '  -1.71%  '

# Row 28: EthereumClassic
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '25.17'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.18%  '

# Row 29: Monero
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '167.18'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.18%  '

# Row 30: Cosmos
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.17'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.24%  '

# Row 31: Toncoin
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.04'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -13.74%  '

# Row 32: InjectiveProtocol
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '33.01'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -5.54%  '

# Row 33: FirstDigitalUSD
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.999'
$ws.Range('D33').ClearFormats()

# Row 34: Celestia
$ws.Range('B34').Value = 'Celestia'
$ws.Range('C34').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.36'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +2.96%  '

# Row 35: Filecoin
$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.12'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +1.52%  '

# Row 36: RenderToken
$ws.Range('E36').Value = '  +0.38%  '

# Row 37: WEMIXToken
$ws.Range('E37').Value = '  -0.50%  '

# Row 38: Hedera
$ws.Range('E38').Value = '  -1.14%  '

# Row 39: Kaspa
$ws.Range('E39').Value = '  +1.03%  '

# Row 41: Stellar
$ws.Range('E41').Value = '  +0.69%  '

# Row 42: LidoDAOToken
$ws.Range('E42').Value = '  -3.00%  '

# Row 43: Maker
$ws.Range('D43').Value = '2.007.60'
$ws.Range('E43').Value = '  +0.46%  '

# Row 44: VeChain
$ws.Range('E44').Value = '  -1.06%  '

# Row 45: ApeXProtocol
$ws.Range('E45').Value = '  -9.02%  '

# Row 46: FraxShare
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.22'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.49%  '

# Row 47: EnergySwap
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '17.44'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.49%  '

# Row 48: NEARProtocol
$ws.Range('E48').Value = '  -2.39%  '

# Row 49: MultiversX
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '54.47'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.38%  '

# Row 50: RocketPoolETH
$ws.Range('D50').Value = '2.536.85'
$ws.Range('E50').Value = '  +0.21%  '

# Row 51: Stacks
$ws.Range('E51').Value = '  -1.46%  '
